# #93 게임서버 session -> character
# Rename "SESSION" -> "CHARACTER" and "MOB | SESSION" -> "MOB | CHARACTER"
# in the OBJECT_TYPE lookup sheet.

$wb = $excel.ActiveWorkbook

# The DESTROY_TYPE tab was previously active/selected; move off of it
# (matches the cursor trail left behind in the saved file).
$wsDestroy = $wb.Worksheets.Item("DESTROY_TYPE")
$wsDestroy.Activate()
$wsDestroy.Range("I11").Select()

$ws = $wb.Worksheets.Item("OBJECT_TYPE")
$ws.Activate()

$ws.Range("A5").Value = "CHARACTER"
$ws.Range("B6").Value = "MOB | CHARACTER"

$ws.Range("B7").Select()
